# This edit re-shuffles the per-row data (Fecha, Volumen, Precio mínimo,
# Precio máximo, Precio promedio ponderado, Unidad de comercialización,
# Origen, Precio $/Kg, Kg / unidad) across the existing data rows (2..24)
# of the weekly Fruta/Hortaliza sheet, while the descriptive columns
# (Mercado ID, Mercado, Región, Codreg, Tipo, Producto ID, Producto,
# Categoría ID, Categoría, Variedad, Calidad) stay identical in every row.
#
# The resulting (target) values for the varying columns of each row are
# written directly as literals below, so the write order does not matter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that vary per row and need to be touched.
# D=4 Fecha, M=13 Volumen, N=14 Precio minimo, O=15 Precio maximo,
# P=16 Precio promedio ponderado, Q=17 Unidad de comercializacion,
# R=18 Origen, S=19 Precio $/Kg, T=20 Kg / unidad

$targets = @(
    @{ Row=2;  D=44432; M=10;  N=35000; O=35000; P=35000; Q='$/caja 18 kilos'; R='Perú'; S=1944; T=18 }
    @{ Row=3;  D=44294; M=15;  N=35000; O=35000; P=35000; Q='$/caja 18 kilos'; R='Región de Arica y Parinacota'; S=1944; T=18 }
    @{ Row=4;  D=44424; M=15;  N=35000; O=35000; P=35000; Q='$/caja 18 kilos'; R='Región de Arica y Parinacota'; S=1944; T=18 }
    @{ Row=5;  D=44379; M=10;  N=30000; O=30000; P=30000; Q='$/caja 18 kilos'; R='Región de Arica y Parinacota'; S=1667; T=18 }
    @{ Row=6;  D=44364; M=90;  N=1700;  O=1700;  P=1700;  Q='$/kilo';          R='Región de Arica y Parinacota'; S=1700; T=1  }
    @{ Row=7;  D=44431; M=30;  N=35000; O=35000; P=35000; Q='$/caja 18 kilos'; R='Región de Arica y Parinacota'; S=1944; T=18 }
    @{ Row=8;  D=44377; M=30;  N=40000; O=40000; P=40000; Q='$/caja 18 kilos'; R='Región de Arica y Parinacota'; S=2222; T=18 }
    @{ Row=9;  D=44435; M=10;  N=35000; O=35000; P=35000; Q='$/caja 18 kilos'; R='Perú'; S=1944; T=18 }
    @{ Row=10; D=44435; M=105; N=35000; O=35000; P=35000; Q='$/caja 18 kilos'; R='Región de Arica y Parinacota'; S=1944; T=18 }
    @{ Row=11; D=44449; M=20;  N=38000; O=38000; P=38000; Q='$/caja 18 kilos'; R='Región de Arica y Parinacota'; S=2111; T=18 }
    @{ Row=12; D=44434; M=40;  N=35000; O=35000; P=35000; Q='$/caja 18 kilos'; R='Región de Arica y Parinacota'; S=1944; T=18 }
    @{ Row=13; D=44405; M=10;  N=35000; O=35000; P=35000; Q='$/caja 18 kilos'; R='Región de Arica y Parinacota'; S=1944; T=18 }
    @{ Row=14; D=44392; M=20;  N=35000; O=35000; P=35000; Q='$/caja 18 kilos'; R='Región de Arica y Parinacota'; S=1944; T=18 }
    @{ Row=15; D=44442; M=15;  N=35000; O=35000; P=35000; Q='$/caja 18 kilos'; R='Perú'; S=1944; T=18 }
    @{ Row=16; D=44448; M=50;  N=38000; O=38000; P=38000; Q='$/caja 18 kilos'; R='Región de Arica y Parinacota'; S=2111; T=18 }
    @{ Row=17; D=44418; M=30;  N=35000; O=35000; P=35000; Q='$/caja 18 kilos'; R='Región de Arica y Parinacota'; S=1944; T=18 }
    @{ Row=18; D=44264; M=20;  N=40000; O=40000; P=40000; Q='$/caja 18 kilos'; R='Región de Arica y Parinacota'; S=2222; T=18 }
    @{ Row=19; D=44369; M=5;   N=35000; O=35000; P=35000; Q='$/caja 18 kilos'; R='Perú'; S=1944; T=18 }
    @{ Row=20; D=44433; M=15;  N=35000; O=35000; P=35000; Q='$/caja 18 kilos'; R='Región de Arica y Parinacota'; S=1944; T=18 }
    @{ Row=21; D=44279; M=30;  N=35000; O=36000; P=35667; Q='$/caja 18 kilos'; R='Región de Arica y Parinacota'; S=1982; T=18 }
    @{ Row=22; D=44363; M=144; N=1700;  O=1700;  P=1700;  Q='$/kilo';          R='Región de Arica y Parinacota'; S=1700; T=1  }
    @{ Row=23; D=44438; M=25;  N=35000; O=35000; P=35000; Q='$/caja 18 kilos'; R='Región de Arica y Parinacota'; S=1944; T=18 }
    @{ Row=24; D=44357; M=10;  N=38000; O=38000; P=38000; Q='$/caja 18 kilos'; R='Perú'; S=2111; T=18 }
)

foreach ($t in $targets) {
    $r = $t.Row
    $ws.Cells.Item($r, 4).Value  = $t.D   # D Fecha
    $ws.Cells.Item($r, 13).Value = $t.M   # M Volumen
    $ws.Cells.Item($r, 14).Value = $t.N   # N Precio minimo
    $ws.Cells.Item($r, 15).Value = $t.O   # O Precio maximo
    $ws.Cells.Item($r, 16).Value = $t.P   # P Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $t.Q   # Q Unidad de comercializacion
    $ws.Cells.Item($r, 18).Value = $t.R   # R Origen
    $ws.Cells.Item($r, 19).Value = $t.S   # S Precio $/Kg
    $ws.Cells.Item($r, 20).Value = $t.T   # T Kg / unidad
}
